$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Description")

# Row 9: "Source from: ABS, Survey of Disability, Ageing and Carers." split into a
# short "Source" label (A9) plus the citation text (B9).
$ws.Range("A9").Value = "Source"
$ws.Range("B9").Value = "ABS, Survey of Disability, Ageing and Carers, various years"

# Row 10: replace the old Northern Territory caveat note with a "References" entry.
$ws.Range("A10").Value = "References"
$ws.Range("B10").Value = "Department of Families, Housing, Community Services and Indigenous Affairs (FaCHSIA), 2009, Shut out: The experience of People with disabilities and their Families in Australia."

# The new reference text uses the workbook's built-in explanatory-text style
# (Calibri 11) with wrapping re-enabled, and the row height is tightened to fit.
$ws.Range("B10").Style = "Excel Built-in Explanatory Text"
$ws.Range("B10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 25.45

$ws.Range("B10").Select()
